# Updates the crypto price/volume table on Sheet1 to reflect the latest
# GitHub Actions scrape (prices, 1h volume %, and a couple of re-ranked rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.225.50'
$ws.Range('E2').Value = '  +1.63%  '
$ws.Range('D3').Value = '1.645.05'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = "'216.91"
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('E6').Value = '  +0.77%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('E9').Value = '  +0.20%  '
$ws.Range('D10').Value = "'19.95"
$ws.Range('E10').Value = '  +1.26%  '
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = "'4.30"
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '1.872.14'
$ws.Range('E13').Value = '  +0.33%  '
$ws.Range('D14').Value = '1.623.96'
$ws.Range('E14').Value = '  -1.03%  '
$ws.Range('E15').Value = '  -2.68%  '
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('D17').Value = "'63.29"
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('D18').Value = '26.223.62'
$ws.Range('E18').Value = '  +1.47%  '
$ws.Range('D20').Value = "'4.45"
$ws.Range('E20').Value = '  -0.80%  '
$ws.Range('D21').Value = "'195.16"
$ws.Range('E21').Value = '  +1.34%  '
$ws.Range('D22').Value = "'10.06"
$ws.Range('E22').Value = '  +0.87%  '
$ws.Range('D23').Value = "'6.32"
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('E24').Value = '  -3.93%  '
$ws.Range('D25').Value = "'143.14"
$ws.Range('E25').Value = '  +0.62%  '
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('E27').Value = '  +1.10%  '
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('D29').Value = "'15.65"
$ws.Range('E29').Value = '  +0.81%  '
$ws.Range('E30').Value = '  +0.67%  '
$ws.Range('E31').Value = '  +1.98%  '
$ws.Range('E32').Value = '  +0.47%  '
$ws.Range('E33').Value = '  +0.58%  '
$ws.Range('E34').Value = '  +1.77%  '
$ws.Range('E35').Value = '  +1.34%  '
$ws.Range('D36').Value = "'0.912"
$ws.Range('E36').Value = '  +0.55%  '
$ws.Range('D37').Value = '1.136.82'
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('D38').Value = "'0.554"
$ws.Range('E38').Value = '  +1.39%  '
$ws.Range('E39').Value = '  -1.76%  '
$ws.Range('E40').Value = '  +1.17%  '
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('D42').Value = "'100.39"
$ws.Range('E42').Value = '  -0.40%  '
$ws.Range('E43').Value = '  -1.30%  '
$ws.Range('D45').Value = '1.780.83'
$ws.Range('E45').Value = '  +0.32%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = "'56.68"
$ws.Range('E46').Value = '  +2.34%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = "'1.48"
$ws.Range('E47').Value = '  +3.71%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = "'0.0517"
$ws.Range('E48').Value = '  +2.96%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = "'0.418"
$ws.Range('E49').Value = '  +0.19%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'7.68"
$ws.Range('E50').Value = '  +2.31%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = "'0.0970"
$ws.Range('E51').Value = '  +1.43%  '

# The apostrophe-prefix trick above stamps a "quote prefix" style onto each
# affected cell. Reset those cells back to the default "Normal" style so the
# saved workbook does not pick up extra style indices that were not present
# in the original file.
$ws.Range('D5').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
